$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.236.25'
$ws.Range("E2").Value = '  +2.36%  '

# Row 3
$ws.Range("D3").Value = '2.596.32'
$ws.Range("E3").Value = '  +1.03%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '533.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.27%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.48%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  +1.72%  '

# Row 9
$ws.Range("D9").Value = '2.613.21'
$ws.Range("E9").Value = '  +1.29%  '

# Row 10
$ws.Range("E10").Value = '  +0.53%  '

# Row 11
$ws.Range("E11").Value = '  +2.98%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.334'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.23%  '

# Row 13
$ws.Range("E13").Value = '  +2.72%  '

# Row 14
$ws.Range("D14").Value = '3.056.96'
$ws.Range("E14").Value = '  +1.08%  '

# Row 15
$ws.Range("D15").Value = '59.183.79'
$ws.Range("E15").Value = '  +2.33%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.28%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.10%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.590.69'
$ws.Range("E18").Value = '  +0.83%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '345.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.45%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.68%  '

# Row 23
$ws.Range("E23").Value = '  -0.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.48%  '

# Row 25
$ws.Range("E25").Value = '  +2.01%  '

# Row 26
$ws.Range("E26").Value = '  +2.13%  '

# Row 27
$ws.Range("E27").Value = '  +0.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.36%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0735'
$ws.Range("E30").Value = '  +3.34%  '

# Row 31
$ws.Range("E31").Value = '  +4.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.88%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.18%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.98'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.36%  '

# Row 36
$ws.Range("E36").Value = '  +1.24%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.24%  '

# Row 38
$ws.Range("E38").Value = '  +4.42%  '

# Row 39
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.846'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.16%  '

# Row 40
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.834'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.35%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.17%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '275.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.22%  '

# Row 44
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.598'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.82%  '

# Row 45
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.71%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0962'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.72%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0521'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.67%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.87%  '

# Row 49
$ws.Range("D49").Value = '1.945.34'
$ws.Range("E49").Value = '  -1.09%  '

# Row 50
$ws.Range("E50").Value = '  +2.43%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.79%  '
